$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.763.41"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.563.48"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'514.28"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "'139.71"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.561"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").Value = "2.577.04"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").Value = "'6.43"
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("D11").Value = "'0.0996"
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("D12").Value = "'0.330"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "'0.132"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "3.015.08"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "57.632.80"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'20.05"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "2.594.27"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "'0.0000132"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").Value = "'332.87"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "'4.25"
$ws.Range("E20").Value = "  -4.19%  "
$ws.Range("D21").Value = "'10.04"
$ws.Range("E21").Value = "  -3.94%  "
$ws.Range("D22").Value = "'6.31"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'64.73"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "'0.396"
$ws.Range("E27").Value = "  -5.13%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.681.64"
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").Value = "'6.91"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "0.0₃0716"
$ws.Range("E31").Value = "  -9.62%  "
$ws.Range("D32").Value = "'6.08"
$ws.Range("E32").Value = "  -5.75%  "
$ws.Range("D33").Value = "'1.56"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "'149.39"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "'18.52"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "'3.91"
$ws.Range("E36").Value = "  -6.16%  "
$ws.Range("D37").Value = "'1.11"
$ws.Range("E37").Value = "  -6.28%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.46"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.835"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'35.74"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "'0.821"
$ws.Range("E41").Value = "  -9.76%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'3.44"
$ws.Range("E43").Value = "  -5.11%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.593"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'10.69"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").Value = "'266.46"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "'0.0941"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "'0.0514"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("D49").Value = "'18.39"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").Value = "1.949.31"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").Value = "'0.0219"
$ws.Range("E51").Value = "  -4.63%  "
